# Update the "data" sheet's time_taken (column F) timestamps to reflect
# the re-run query time.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("data")

$ws.Range("F2").Value = "2021-10-05 14:20:21.858948"
$ws.Range("F3").Value = "2021-10-05 14:20:21.858957"
$ws.Range("F4").Value = "2021-10-05 14:20:21.858960"
$ws.Range("F5").Value = "2021-10-05 14:20:21.858963"
$ws.Range("F6").Value = "2021-10-05 14:20:21.858967"
$ws.Range("F7").Value = "2021-10-05 14:20:21.858970"
$ws.Range("F8").Value = "2021-10-05 14:20:21.858972"
$ws.Range("F9").Value = "2021-10-05 14:20:21.858975"
$ws.Range("F10").Value = "2021-10-05 14:20:21.858978"
$ws.Range("F11").Value = "2021-10-05 14:20:21.858981"
$ws.Range("F12").Value = "2021-10-05 14:20:21.858984"
$ws.Range("F13").Value = "2021-10-05 14:20:21.858987"
$ws.Range("F14").Value = "2021-10-05 14:20:21.858990"
$ws.Range("F15").Value = "2021-10-05 14:20:21.858993"
$ws.Range("F16").Value = "2021-10-05 14:20:21.858996"
$ws.Range("F17").Value = "2021-10-05 14:20:21.858998"
$ws.Range("F18").Value = "2021-10-05 14:20:21.859001"
$ws.Range("F19").Value = "2021-10-05 14:20:21.859004"

# Add a new "metadata" worksheet after "data", describing the panel
# query that produced the data sheet.
$meta = $wb.Worksheets.Add($null, $ws)
$meta.Name = "metadata"

$meta.Range("B1").Value = "data_name"
$meta.Range("C1").Value = "data_id"
$meta.Range("D1").Value = "data_version"
$meta.Range("E1").Value = "data_version_created"
$meta.Range("F1").Value = "panel_query_time"
$meta.Range("G1").Value = "panel_get_request"

# Match the bold/centered header style used on row 1 of the "data" sheet.
$ws.Range("B1").Copy()
$meta.Range("B1:G1").PasteSpecial(-4122)  # xlPasteFormats

$meta.Range("A2").Value = 0
$ws.Range("A2").Copy()
$meta.Range("A2").PasteSpecial(-4122)  # xlPasteFormats
$meta.Range("B2").Value = "Familial rhabdomyosarcoma"
$meta.Range("C2").Value = 290

# data_version ("1.4") must stay a text value, not be coerced to a number.
$verCell = $meta.Range("D2")
$verCell.NumberFormat = "@"
$verCell.Value = "1.4"
$verCell.Style = "Normal"

$meta.Range("E2").Value = "2019-01-29T15:49:45.308605Z"
$meta.Range("F2").Value = "2021-10-05 14:20:21.855393"
$meta.Range("G2").Value = "https://panelapp.genomicsengland.co.uk/api/v1/panels/290/?format=json"

Write-Output "Edit complete"
